{"js": "// Change the label \"CRSTb\" to \"CRESET\" (device pin naming fix).\nconst results = context.document.body.search(\"CRSTb\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"CRESET\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Change the label \"CRSTb\" to \"CRESET\" (device pin naming fix).\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"CRSTb\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"CRESET\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n"}
